$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 218 (A218 = "01-01-2021") ---
$ws.Range("B218").Value = 36657
$ws.Range("E218").Value = 36498
$ws.Range("G218").Value = 1062

# --- Row 219 (A219 = "01-02-2021") ---
$ws.Range("B219").Value = 35564
$ws.Range("E219").Value = 35402
$ws.Range("G219").Value = 1063

# --- Row 220 (A220 = "01-03-2021") ---
$ws.Range("B220").Value = 33754
$ws.Range("E220").Value = 33474
$ws.Range("G220").Value = 1050

# --- Row 221 (A221 = "01-04-2021") ---
$ws.Range("B221").Value = 35862
$ws.Range("G221").Value = 1044

# --- Row 222 (A222 = "01-05-2021") ---
$ws.Range("B222").Value = 37594
$ws.Range("G222").Value = 1044

# --- New row 223 (A223 = "01-06-2021") ---
# Build the date-like label through a formula first, then collapse it to a
# plain cached value via copy/paste-values, so Excel's smart "looks like a
# date" literal parser (which would otherwise convert "01-06-2021" into a
# date serial and stamp a new number-format style on the cell) never gets a
# chance to fire.
$ws.Range("A223").Formula = '="01-06-2021"'
$ws.Range("A223").Copy()
$ws.Range("A223").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B223").Value = 37860
$ws.Range("C223").Value = 142
$ws.Range("D223").Value = 142
$ws.Range("E223").Value = 37718
$ws.Range("F223").Value = 36692
$ws.Range("G223").Value = 1026
